$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2 = "2025-08-13T01:55:53.253319+09:00"
    3 = "2025-08-13T01:55:53.285291+09:00"
    4 = "2025-08-13T01:55:53.321593+09:00"
    5 = "2025-08-13T01:55:53.367439+09:00"
    6 = "2025-08-13T01:55:53.423481+09:00"
    7 = "2025-08-13T01:55:53.482434+09:00"
    8 = "2025-08-13T01:55:53.538615+09:00"
    9 = "2025-08-13T01:55:53.585426+09:00"
    10 = "2025-08-13T01:55:53.625601+09:00"
    11 = "2025-08-13T01:55:53.665736+09:00"
    12 = "2025-08-13T01:55:53.685832+09:00"
    13 = "2025-08-13T01:55:53.704012+09:00"
    14 = "2025-08-13T01:55:53.722664+09:00"
    15 = "2025-08-13T01:55:53.767675+09:00"
    16 = "2025-08-13T01:55:53.823886+09:00"
    17 = "2025-08-13T01:55:53.859082+09:00"
    18 = "2025-08-13T01:55:53.909603+09:00"
    19 = "2025-08-13T01:55:53.950572+09:00"
    20 = "2025-08-13T01:55:53.995393+09:00"
    21 = "2025-08-13T01:55:54.047876+09:00"
    22 = "2025-08-13T01:55:54.102602+09:00"
    23 = "2025-08-13T01:55:54.145470+09:00"
    24 = "2025-08-13T01:55:54.190297+09:00"
    25 = "2025-08-13T01:55:54.234260+09:00"
    26 = "2025-08-13T01:55:54.274789+09:00"
    27 = "2025-08-13T01:55:54.314644+09:00"
    28 = "2025-08-13T01:55:54.366434+09:00"
    29 = "2025-08-13T01:55:54.422667+09:00"
    30 = "2025-08-13T01:55:54.477459+09:00"
    31 = "2025-08-13T01:55:54.524627+09:00"
    32 = "2025-08-13T01:55:54.574708+09:00"
    33 = "2025-08-13T01:55:54.630188+09:00"
    34 = "2025-08-13T01:55:54.676484+09:00"
    35 = "2025-08-13T01:55:54.701121+09:00"
    36 = "2025-08-13T01:55:54.722034+09:00"
    37 = "2025-08-13T01:55:54.751314+09:00"
    38 = "2025-08-13T01:55:54.786823+09:00"
    39 = "2025-08-13T01:55:54.837621+09:00"
    40 = "2025-08-13T01:55:54.888795+09:00"
    41 = "2025-08-13T01:55:54.928092+09:00"
    42 = "2025-08-13T01:55:54.974682+09:00"
    43 = "2025-08-13T01:55:55.019904+09:00"
    44 = "2025-08-13T01:55:55.053073+09:00"
    45 = "2025-08-13T01:55:55.091615+09:00"
    46 = "2025-08-13T01:55:55.133436+09:00"
    47 = "2025-08-13T01:55:55.176064+09:00"
    48 = "2025-08-13T01:55:55.218868+09:00"
    49 = "2025-08-13T01:55:55.255322+09:00"
    50 = "2025-08-13T01:55:55.296021+09:00"
    51 = "2025-08-13T01:55:55.330133+09:00"
    52 = "2025-08-13T01:55:55.375853+09:00"
    53 = "2025-08-13T01:55:55.415729+09:00"
    54 = "2025-08-13T01:55:55.458289+09:00"
    55 = "2025-08-13T01:55:55.500614+09:00"
    56 = "2025-08-13T01:55:55.544531+09:00"
    57 = "2025-08-13T01:55:55.586892+09:00"
    58 = "2025-08-13T01:55:55.640140+09:00"
    59 = "2025-08-13T01:55:55.692077+09:00"
    60 = "2025-08-13T01:55:55.729388+09:00"
    61 = "2025-08-13T01:55:55.773001+09:00"
    62 = "2025-08-13T01:55:55.819310+09:00"
    63 = "2025-08-13T01:55:55.862717+09:00"
    64 = "2025-08-13T01:55:55.910090+09:00"
    65 = "2025-08-13T01:55:55.952182+09:00"
    66 = "2025-08-13T01:55:55.994977+09:00"
    67 = "2025-08-13T01:55:56.040061+09:00"
    68 = "2025-08-13T01:55:56.090615+09:00"
    69 = "2025-08-13T01:55:56.138924+09:00"
    70 = "2025-08-13T01:55:56.188462+09:00"
    71 = "2025-08-13T01:55:56.238163+09:00"
    72 = "2025-08-13T01:55:56.278651+09:00"
    73 = "2025-08-13T01:55:56.309434+09:00"
    74 = "2025-08-13T01:55:56.336357+09:00"
    75 = "2025-08-13T01:55:56.377485+09:00"
    76 = "2025-08-13T01:55:56.429566+09:00"
    77 = "2025-08-13T01:55:56.478229+09:00"
    78 = "2025-08-13T01:55:56.514014+09:00"
    79 = "2025-08-13T01:55:56.566186+09:00"
    80 = "2025-08-13T01:55:56.618902+09:00"
    81 = "2025-08-13T01:55:56.667213+09:00"
    82 = "2025-08-13T01:55:56.711279+09:00"
    83 = "2025-08-13T01:55:56.751206+09:00"
    84 = "2025-08-13T01:55:56.795307+09:00"
    85 = "2025-08-13T01:55:56.835927+09:00"
    86 = "2025-08-13T01:55:56.875197+09:00"
    87 = "2025-08-13T01:55:56.922872+09:00"
    88 = "2025-08-13T01:55:56.966119+09:00"
    89 = "2025-08-13T01:55:57.017825+09:00"
    90 = "2025-08-13T01:55:57.060157+09:00"
    91 = "2025-08-13T01:55:57.092958+09:00"
    92 = "2025-08-13T01:55:57.122578+09:00"
    93 = "2025-08-13T01:55:57.145922+09:00"
    94 = "2025-08-13T01:55:57.192625+09:00"
    95 = "2025-08-13T01:55:57.240699+09:00"
    96 = "2025-08-13T01:55:57.293206+09:00"
    97 = "2025-08-13T01:55:57.337886+09:00"
    98 = "2025-08-13T01:55:57.385130+09:00"
    99 = "2025-08-13T01:55:57.436720+09:00"
    100 = "2025-08-13T01:55:57.493456+09:00"
    101 = "2025-08-13T01:55:57.531219+09:00"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $timestamps[$row]
}
